$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-127 shift down to 11-128,
# carrying their formatting/styles with them (matches the diff exactly).
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new record.
$ws.Cells.Item(10,1).Value = 7
$ws.Cells.Item(10,2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(10,3).Value = 'Ñuble'
$ws.Cells.Item(10,4).Value = 44630
$ws.Cells.Item(10,5).Value = 16
$ws.Cells.Item(10,6).Value = 'Fruta'
$ws.Cells.Item(10,7).Value = 100109
$ws.Cells.Item(10,8).Value = 'Uva'
$ws.Cells.Item(10,9).Value = 100109001
$ws.Cells.Item(10,10).Value = 'Uva'
$ws.Cells.Item(10,11).Value = 'Red Globe'
$ws.Cells.Item(10,12).Value = 'Primera'
$ws.Cells.Item(10,13).Value = 120
$ws.Cells.Item(10,14).Value = 10000
$ws.Cells.Item(10,15).Value = 11000
$ws.Cells.Item(10,16).Value = 10500
$ws.Cells.Item(10,17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(10,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10,19).Value = 583
$ws.Cells.Item(10,20).Value = 18
